# Laborator 19.03.2024 - checked "saptamana 4" (column F) presence for the
# students that attended, added 3 new students, then re-sorted the roster
# alphabetically by name (same way the sheet was already sorted).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Mark "saptamana 4" (column F) attendance for students that were present ---
$week4Names = @(
    "Adina Pop",
    "Alexandru Lupse",
    "Andreea Farcas",
    "Aniko Vieriu",
    "Bianca Abrudan",
    "Cristina Nemcea",
    "Daria Petre",
    "David Florea",
    "David Klein",
    "Erik Lazin",
    "Nuria Girz",
    "Patrick Tocut",
    "Razvan Ardeli",
    "Renata Halasz",
    "Roland Roman",
    "Vanesa Crepce",
    "Vlad Chis",
    "Raluca Veres",
    "Gabriela Maghear"
)

for ($r = 3; $r -le 34; $r++) {
    $name = $ws.Cells.Item($r, 2).Value2
    if ($week4Names -contains $name) {
        $ws.Cells.Item($r, 6).Value = $true
    }
}

# --- 2. Add the 3 new students, present only for "saptamana 4" ---
$ws.Range("B35").Value = "Luminita Hava"
$ws.Range("F35").Value = $true

$ws.Range("B36").Value = "Victor Balaj"
$ws.Range("F36").Value = $true

$ws.Range("B37").Value = "Emanuel Socaciu"
$ws.Range("F37").Value = $true

# --- 3. Re-sort the whole roster (B3:S37) alphabetically by name, like before ---
$sortRange = $ws.Range("B3:S37")
$sortKey = $ws.Range("B3")
$sortRange.Sort($sortKey)

# --- 4. Update the view: clear the old frozen top row and move the selection ---
$ws.Range("N12").Select()
